# feat: remove export pdf and remove all nilai column
#
# Clears the "nilai1".."nilai16" header cells (J1:Y1) while keeping their
# existing cell style, and moves the active selection to P19.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the content of the nilai1..nilai16 header cells (J1:Y1), leaving the
# cells present (with their style) but empty.
$ws.Range("J1:Y1").ClearContents()

# Update the saved selection to match the target workbook.
$ws.Range("P19").Select()
